$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "1.070" "1.016"
Replace-Text "0.996" "0.999"
Replace-Text "0.282" "0.346"
Replace-Text "3.061" "3.147"
Replace-Text "0.939" "0.962"
Replace-Text "0.741" "1.149"
Replace-Text "0.597" "0.756"
Replace-Text "7.134" "7.469"
Replace-Text "0.781" "0.858"
